$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 14 (for date 2022-04-25 / serial 44676),
# pushing all the existing rows from 14..77 down to 15..78.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 6
$ws.Range("B14").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 44676
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100114007
$ws.Range("G14").Value = "Jengibre"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 270
$ws.Range("K14").Value = 11000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11481
$ws.Range("N14").Value = "`$/caja 13 kilos"
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 883
$ws.Range("Q14").Value = 13
$ws.Range("R14").Value = "Hortaliza"
